$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cutting Speed")

$data = @(
    @(120, 40, 4, 0.12),
    @(90, 8, 4, 0.1),
    @(60, 8, 4, 0.1),
    @(40, 8, 4, 0.1),
    @(40, 10, 4, 0.1),
    @(40, 10, 4, 0.01),
    @(30, 10, 4, 0.01),
    @(60, 10, 4, 0.01),
    @(60, 12, 4, 0.01),
    @(210, 80, 6, 0.12),
    @(300, 80, 6, 0.12),
    @(300, 60, 6, 0.12),
    @(300, 60, 4, 0.12),
    @(300, 60, 4, 0.1)
)

$startRow = 30
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
}
